$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1923076923076923
$ws.Range("C2").Value = 0.5416666666666666
$ws.Range("J2").Value = 0.02564102564102564
$ws.Range("P2").Value = 0.1282051282051282
$ws.Range("S2").Value = 0.1121794871794872
$ws.Range("B3").Value = 0.005952380952380952
$ws.Range("C3").Value = 0.0119047619047619
$ws.Range("J3").Value = 0.02976190476190476
$ws.Range("P3").Value = 0.7440476190476191
$ws.Range("S3").Value = 0.2083333333333333
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.04471544715447155
$ws.Range("D6").Value = 0.01626016260162602
$ws.Range("F6").Value = 0.04471544715447155
$ws.Range("J6").Value = 0.2682926829268293
$ws.Range("O6").Value = 0.04471544715447155
$ws.Range("Q6").Value = 0.2154471544715447
$ws.Range("R6").Value = 0.04878048780487805
$ws.Range("S6").Value = 0.3170731707317073
$ws.Range("B7").Value = 0.1308016877637131
$ws.Range("D7").Value = 0.01687763713080169
$ws.Range("F7").Value = 0.03375527426160337
$ws.Range("J7").Value = 0.1265822784810127
$ws.Range("O7").Value = 0.02109704641350211
$ws.Range("Q7").Value = 0.1687763713080169
$ws.Range("R7").Value = 0.0759493670886076
$ws.Range("S7").Value = 0.4261603375527426
$ws.Range("B8").Value = 0.08108108108108109
$ws.Range("D8").Value = 0.02316602316602316
$ws.Range("E8").Value = 0.001930501930501931
$ws.Range("F8").Value = 0.05019305019305019
$ws.Range("J8").Value = 0.138996138996139
$ws.Range("O8").Value = 0.03088803088803089
$ws.Range("Q8").Value = 0.1988416988416988
$ws.Range("R8").Value = 0.07528957528957529
$ws.Range("S8").Value = 0.3996138996138996
$ws.Range("B9").Value = 0.0778688524590164
$ws.Range("D9").Value = 0.02459016393442623
$ws.Range("F9").Value = 0.06557377049180328
$ws.Range("J9").Value = 0.0778688524590164
$ws.Range("O9").Value = 0.01229508196721311
$ws.Range("Q9").Value = 0.2131147540983606
$ws.Range("R9").Value = 0.1229508196721311
$ws.Range("S9").Value = 0.4057377049180328
$ws.Range("B10").Value = 0.1122529644268775
$ws.Range("D10").Value = 0.01660079051383399
$ws.Range("E10").Value = 0.0007905138339920949
$ws.Range("F10").Value = 0.09486166007905138
$ws.Range("J10").Value = 0.1185770750988142
$ws.Range("O10").Value = 0.01818181818181818
$ws.Range("Q10").Value = 0.1936758893280632
$ws.Range("R10").Value = 0.09090909090909091
$ws.Range("S10").Value = 0.3541501976284585
$ws.Range("G11").Value = 0.1362229102167183
$ws.Range("J11").Value = 0.0804953560371517
$ws.Range("K11").Value = 0.1671826625386997
$ws.Range("L11").Value = 0.6068111455108359
$ws.Range("S11").Value = 0.009287925696594427
$ws.Range("G12").Value = 0.7725118483412322
$ws.Range("J12").Value = 0.1611374407582938
$ws.Range("K12").Value = 0.01421800947867299
$ws.Range("L12").Value = 0.04265402843601896
$ws.Range("S12").Value = 0.009478672985781991
$ws.Range("G13").Value = 0.6842105263157895
$ws.Range("J13").Value = 0.2631578947368421
$ws.Range("S13").Value = 0.05263157894736842
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.0371900826446281
$ws.Range("H15").Value = 0.1652892561983471
$ws.Range("I15").Value = 0.06611570247933884
$ws.Range("J15").Value = 0.2933884297520661
$ws.Range("K15").Value = 0.07024793388429752
$ws.Range("M15").Value = 0.01652892561983471
$ws.Range("O15").Value = 0.04132231404958678
$ws.Range("S15").Value = 0.3099173553719008
$ws.Range("F16").Value = 0.01047120418848168
$ws.Range("H16").Value = 0.2041884816753927
$ws.Range("I16").Value = 0.1047120418848168
$ws.Range("J16").Value = 0.3141361256544503
$ws.Range("K16").Value = 0.1047120418848168
$ws.Range("M16").Value = 0.02617801047120419
$ws.Range("O16").Value = 0.06282722513089005
$ws.Range("S16").Value = 0.1727748691099476
$ws.Range("F17").Value = 0.02240325865580448
$ws.Range("H17").Value = 0.2016293279022403
$ws.Range("I17").Value = 0.120162932790224
$ws.Range("J17").Value = 0.3604887983706721
$ws.Range("K17").Value = 0.1140529531568228
$ws.Range("M17").Value = 0.02036659877800407
$ws.Range("N17").Value = 0.002036659877800407
$ws.Range("O17").Value = 0.05295315682281059
$ws.Range("S17").Value = 0.1059063136456212
$ws.Range("F18").Value = 0.009345794392523364
$ws.Range("H18").Value = 0.1682242990654206
$ws.Range("I18").Value = 0.1074766355140187
$ws.Range("J18").Value = 0.3925233644859813
$ws.Range("K18").Value = 0.09813084112149532
$ws.Range("M18").Value = 0.009345794392523364
$ws.Range("O18").Value = 0.07009345794392523
$ws.Range("S18").Value = 0.1448598130841121
$ws.Range("F19").Value = 0.01239970824215901
$ws.Range("H19").Value = 0.2261123267687819
$ws.Range("I19").Value = 0.09336250911743253
$ws.Range("J19").Value = 0.3391684901531729
$ws.Range("K19").Value = 0.1057622173595915
$ws.Range("M19").Value = 0.02771699489423778
$ws.Range("N19").Value = 0.001458789204959883
$ws.Range("O19").Value = 0.06710430342815463
$ws.Range("S19").Value = 0.1269146608315098
